$d = $word.ActiveDocument

# Replace the title text "What type of App" with "Describe what we did"
# (Find/Replace also collapses the two runs into a single run.)
$d.Content.Find.Execute("What type of App", $false, $false, $false, $false, $false, $true, 1, $false, "Describe what we did", 2)

# Remove the paragraphs between the title paragraph and the bookmark
# paragraph (Problem:/Solution:/Data Analysis:/Data source: and the
# blank paragraphs interleaved with them).
$r1 = $d.Range($d.Paragraphs.Item(2).Range.Start, $d.Paragraphs.Item(9).Range.End)
$r1.Delete()

# Remove every paragraph after the bookmark paragraph (Model used:,
# Easy, Linear regression, the indented code-like lines, Hard, and the
# closing paragraph) while preserving the final sectPr.
$count = $d.Paragraphs.Count
$r2 = $d.Range($d.Paragraphs.Item(3).Range.Start, $d.Paragraphs.Item($count).Range.End)
$r2.Delete()
